# "3.1 Activity-PSS_ Turtle.pptx" edit:
# Remove the pre-drawn answer-key shapes (square + coordinate labels) that
# were overlaid on the blank turtle-graphics grids of the "Problem 1" and
# "Problem 2" slides, leaving the grids blank again.
#
# Problem 1 slide: shapes 471-475 (Google Shape;471;p41 .. Google Shape;475;p41)
# Problem 2 slide: shapes 494-498 (Google Shape;494;p42 .. Google Shape;498;p42)

$p = $ppt.ActivePresentation

$idsToRemove = @(471, 472, 473, 474, 475, 494, 495, 496, 497, 498)

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = $slide.Shapes.Count; $shi -ge 1; $shi--) {
        $shape = $slide.Shapes.Item($shi)
        if ($idsToRemove -contains $shape.Id) {
            $shape.Delete()
        }
    }
}
